{"js": "// Apply the Golden Tiger content refresh:\n// - New title/headline copy (appears twice: H1 and bold summary line)\n// - Pros list: rewording of 2 bullets, \"structure\" -> \"design\"\n// - Cons list: rewording of 2 bullets\n// - Meta description rewrite (italic summary line)\n\nconst replacements = [\n  [\n    \"Play Golden Tiger Slot Free - Simple Gameplay & Big Wins\",\n    \"Play Golden Tiger - Free Online Slot Game\",\n  ],\n  [\n    \"Simple and minimalistic structure\",\n    \"Simple and minimalistic design\",\n  ],\n  [\n    \"Double bonus wheel function offered\",\n    \"Straightforward winning potential\",\n  ],\n  [\n    \"Good RTP of 95.95%\",\n    \"Good Return to Player (RTP)\",\n  ],\n  [\n    \"No Wilds, Scatters or special bonuses\",\n    \"Lack of Wilds, Scatters, and bonus features\",\n  ],\n  [\n    \"May be too simple for some players\",\n    \"May be too simple for players seeking advanced gameplay\",\n  ],\n  [\n    \"Play Golden Tiger slot online for free! Simple gameplay with 5 pay lines, bonus game, and RTP of 95.95%. Ideal for new players or those after a relaxed gaming experience.\",\n    \"Play Golden Tiger for free and enjoy a simple and straightforward gaming experience.\",\n  ],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the Golden Tiger content refresh:\n# - New title/headline copy (appears twice: H1 and bold summary line)\n# - Pros list: rewording of 2 bullets, \"structure\" -> \"design\"\n# - Cons list: rewording of 2 bullets\n# - Meta description rewrite (italic summary line)\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Play Golden Tiger Slot Free - Simple Gameplay & Big Wins\", \"Play Golden Tiger - Free Online Slot Game\"),\n    @(\"Simple and minimalistic structure\", \"Simple and minimalistic design\"),\n    @(\"Double bonus wheel function offered\", \"Straightforward winning potential\"),\n    @(\"Good RTP of 95.95%\", \"Good Return to Player (RTP)\"),\n    @(\"No Wilds, Scatters or special bonuses\", \"Lack of Wilds, Scatters, and bonus features\"),\n    @(\"May be too simple for some players\", \"May be too simple for players seeking advanced gameplay\"),\n    @(\"Play Golden Tiger slot online for free! Simple gameplay with 5 pay lines, bonus game, and RTP of 95.95%. Ideal for new players or those after a relaxed gaming experience.\", \"Play Golden Tiger for free and enjoy a simple and straightforward gaming experience.\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
